$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2,4,5,6,7) have effectively been rotated/shuffled between
# each other (row 3 is untouched). Apply the new values directly per the
# target state described by the diff.

$rowData = @{
    2 = @{ D = 44382; J = 160; K = 7000;  L = 8000;  M = 7438; P = 124 }
    4 = @{ D = 44281; J = 120; K = 5500;  L = 6000;  M = 5750; P = 96  }
    5 = @{ D = 44362; J = 120; K = 8000;  L = 9000;  M = 8500; P = 142 }
    6 = @{ D = 44421; J = 100; K = 8000;  L = 9000;  M = 8500; P = 142 }
    7 = @{ D = 44400; J = 120; K = 9000;  L = 10000; M = 9500; P = 158 }
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
